$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01297332360334063
$ws.Range("C2").Value = 0.01508271164450211
$ws.Range("D2").Value = 0.01349180854481208
$ws.Range("E2").Value = 0.01378974691758598
$ws.Range("F2").Value = 0.01546853374582757
$ws.Range("G2").Value = 0.01357246728065923
$ws.Range("H2").Value = 0.01337982484592929
$ws.Range("I2").Value = 0.01368088723387031
$ws.Range("J2").Value = 0.0130048465266559
$ws.Range("K2").Value = 0.01418897129958078
$ws.Range("M2").Value = 0.01422797089733225
$ws.Range("N2").Value = 0.01384615384615385
$ws.Range("O2").Value = 0.01295756397797214
$ws.Range("P2").Value = 0.01392
$ws.Range("Q2").Value = 0.01376811594202899
$ws.Range("R2").Value = 0.01250710631040364
$ws.Range("S2").Value = 0.01257198475139914
$ws.Range("T2").Value = 0.01325895383620341
$ws.Range("U2").Value = 0.01322473994032739
$ws.Range("V2").Value = 0.01345580533397792
$ws.Range("W2").Value = 0.01307241960060951
$ws.Range("X2").Value = 0.01499032882011605
$ws.Range("Y2").Value = 0.01436990393154113
$ws.Range("Z2").Value = 0.01511826383808827
$ws.Range("AA2").Value = 0.01425499721271004
$ws.Range("AB2").Value = 0.01456623209399646
$ws.Range("AC2").Value = 0.01397189468583428
$ws.Range("AD2").Value = 0.01417233560090703
$ws.Range("AE2").Value = 0.01401982112642011
$ws.Range("AF2").Value = 0.01285223916340704
$ws.Range("AG2").Value = 0.01319176319176319
$ws.Range("AH2").Value = 0.01392168037124481
$ws.Range("AI2").Value = 0.01384839650145773
$ws.Range("AJ2").Value = 0.0149640055002831
$ws.Range("AK2").Value = 0.01311027105284324
$ws.Range("AL2").Value = 0.01402886398451987
$ws.Range("AM2").Value = 0.01447410743004181
$ws.Range("AN2").Value = 0.01398151868220169
$ws.Range("AO2").Value = 0.01323745710083347
$ws.Range("AP2").Value = 0.01391871470611628
$ws.Range("AQ2").Value = 0.01277877258668403
$ws.Range("AR2").Value = 0.01354620222544751
$ws.Range("AS2").Value = 0.01301727601446364
$ws.Range("AT2").Value = 0.01247165532879819
$ws.Range("AU2").Value = 0.01445274439753167
$ws.Range("AV2").Value = 0.01457725947521866
$ws.Range("AW2").Value = 0.01417639951671365
$ws.Range("AX2").Value = 0.01259241205621903
$ws.Range("AY2").Value = 0.01345146999597261
$ws.Range("AZ2").Value = 0.01376073480419796
$ws.Range("BA2").Value = 0.0007317412049066707
$ws.Range("B3").Value = 0.272108843537415
$ws.Range("C3").Value = 0.303921568627451
$ws.Range("D3").Value = 0.2926829268292683
$ws.Range("E3").Value = 0.2891156462585034
$ws.Range("F3").Value = 0.3059581320450886
$ws.Range("G3").Value = 0.2828282828282828
$ws.Range("H3").Value = 0.2665589660743134
$ws.Range("I3").Value = 0.2793388429752066
$ws.Range("J3").Value = 0.2771084337349398
$ws.Range("K3").Value = 0.2953020134228188
$ws.Range("M3").Value = 0.2993197278911565
$ws.Range("N3").Value = 0.283112582781457
$ws.Range("O3").Value = 0.2576489533011272
$ws.Range("P3").Value = 0.2929292929292929
$ws.Range("Q3").Value = 0.2817133443163097
$ws.Range("R3").Value = 0.2592592592592592
$ws.Range("S3").Value = 0.2566225165562914
$ws.Range("T3").Value = 0.2701812191103789
$ws.Range("U3").Value = 0.2789115646258503
$ws.Range("V3").Value = 0.2742200328407225
$ws.Range("W3").Value = 0.2810344827586207
$ws.Range("X3").Value = 0.3079470198675497
$ws.Range("Y3").Value = 0.2899022801302932
$ws.Range("Z3").Value = 0.2980769230769231
$ws.Range("AA3").Value = 0.2905844155844156
$ws.Range("AB3").Value = 0.3083475298126065
$ws.Range("AC3").Value = 0.2836065573770492
$ws.Range("AD3").Value = 0.2887788778877888
$ws.Range("AE3").Value = 0.2890365448504983
$ws.Range("AF3").Value = 0.2648026315789473
$ws.Range("AG3").Value = 0.277027027027027
$ws.Range("AH3").Value = 0.2785016286644951
$ws.Range("AI3").Value = 0.280327868852459
$ws.Range("AJ3").Value = 0.2988691437802908
$ws.Range("AK3").Value = 0.2583201267828843
$ws.Range("AL3").Value = 0.2909698996655518
$ws.Range("AM3").Value = 0.2907915993537964
$ws.Range("AN3").Value = 0.2810985460420032
$ws.Range("AO3").Value = 0.2718120805369127
$ws.Range("AP3").Value = 0.2966101694915254
$ws.Range("AQ3").Value = 0.2612312811980033
$ws.Range("AR3").Value = 0.28
$ws.Range("AS3").Value = 0.2686567164179104
$ws.Range("AT3").Value = 0.2601351351351351
$ws.Range("AU3").Value = 0.2894308943089431
$ws.Range("AV3").Value = 0.2926829268292683
$ws.Range("AW3").Value = 0.2811501597444089
$ws.Range("AX3").Value = 0.2587646076794658
$ws.Range("AY3").Value = 0.2840136054421769
$ws.Range("AZ3").Value = 0.2822563641648795
$ws.Range("BA3").Value = 0.01382624494248351
$ws.Range("B4").Value = 0.02476588499342156
$ws.Range("C4").Value = 0.02873918417799753
$ws.Range("D4").Value = 0.0257945647167204
$ws.Range("E4").Value = 0.02632393930009291
$ws.Range("F4").Value = 0.02944823310601364
$ws.Range("G4").Value = 0.02590194264569843
$ws.Range("H4").Value = 0.02548065786425759
$ws.Range("I4").Value = 0.02608427226423831
$ws.Range("J4").Value = 0.02484376205539696
$ws.Range("K4").Value = 0.02707692307692307
$ws.Range("M4").Value = 0.02716468590831918
$ws.Range("N4").Value = 0.02640111162575266
$ws.Range("O4").Value = 0.02467422314750559
$ws.Range("P4").Value = 0.02657705819459294
$ws.Range("Q4").Value = 0.02625316650034544
$ws.Range("R4").Value = 0.02386302006663051
$ws.Range("S4").Value = 0.02396968994046238
$ws.Range("T4").Value = 0.02527743526510481
$ws.Range("U4").Value = 0.02525213642312726
$ws.Range("V4").Value = 0.02565284178187404
$ws.Range("W4").Value = 0.02498275729940992
$ws.Range("X4").Value = 0.02858899477405472
$ws.Range("Y4").Value = 0.02738250903776632
$ws.Range("Z4").Value = 0.02877697841726619
$ws.Range("AA4").Value = 0.02717680103241478
$ws.Range("AB4").Value = 0.0278183355106432
$ws.Range("AC4").Value = 0.02663177339901478
$ws.Range("AD4").Value = 0.02701868148834337
$ws.Range("AE4").Value = 0.02674248828094982
$ws.Range("AF4").Value = 0.02451465550057099
$ws.Range("AG4").Value = 0.02518427518427519
$ws.Range("AH4").Value = 0.02651779483600837
$ws.Range("AI4").Value = 0.02639296187683284
$ws.Range("AJ4").Value = 0.0285010013865352
$ws.Range("AK4").Value = 0.02495407225964483
$ws.Range("AL4").Value = 0.02676717175601877
$ws.Range("AM4").Value = 0.02757564151666029
$ws.Range("AN4").Value = 0.02663808940600122
$ws.Range("AO4").Value = 0.02524544179523142
$ws.Range("AP4").Value = 0.02658968320291727
$ws.Range("AQ4").Value = 0.02436563979203849
$ws.Range("AR4").Value = 0.02584217812644209
$ws.Range("AS4").Value = 0.02483139178418148
$ws.Range("AT4").Value = 0.02380216383307574
$ws.Range("AU4").Value = 0.02753074008197355
$ws.Range("AV4").Value = 0.0277713492247165
$ws.Range("AW4").Value = 0.02699179510773713
$ws.Range("AX4").Value = 0.02401611403780601
$ws.Range("AY4").Value = 0.02568638006613858
$ws.Range("AZ4").Value = 0.02624146506416725
$ws.Range("BA4").Value = 0.001384356003686979
